$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Supplier row (row 9): Confirmed/Imported counts increased 1366 -> 1423
$ws.Range("E9").Value = 1423
$ws.Range("F9").Value = 1423

# Update Product row (row 11): values were previously empty, now populated
$ws.Range("E11").Value = 4528
$ws.Range("F11").Value = 4528

# Re-assert the H11 formula so its cached result refreshes from the stale
# empty-string ("") state to the numeric 1 now that E11/F11 are populated
$ws.Range("H11").Formula = "=IF(E11>0, F11/E11, """")"

# Move the active selection to F12
$ws.Range("F12").Select()
